# Regenerate the "K" column (column G) values for save_data sheet.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the observable effect on this sheet is that the
# values in column G (header "K") are recalculated/rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 2
    12 = 1
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 2
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 1
    37 = 1
    38 = 0
    39 = 1
    40 = 0
    41 = 2
    43 = 3
    44 = 1
    45 = 2
    46 = 1
    47 = 2
    48 = 0
    49 = 1
    50 = 2
    51 = 1
    52 = 2
    53 = 0
    54 = 0
    55 = 3
    56 = 3
    58 = 1
    59 = 2
    60 = 0
    61 = 1
    62 = 2
    63 = 4
    64 = 3
    65 = 5
    66 = 0
    67 = 1
    68 = 1
    69 = 0
    70 = 2
    71 = 0
    72 = 0
    73 = 1
    74 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
